$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "CasesTab" row (row 2, column B): this query previously lived in the
# StudyFilesTab row; swap it into the CasesTab row.
$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis), (r:registration)-->(c), (samp:sample)-->(c)
WHERE s.clinical_study_designation IN ['OSA01'] and demo.breed in ['Unknown'] and samp.summarized_sample_type IN ['Normal Cell Line']

OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co,demo.patient_age_at_enrollment AS age, demo.weight as weight
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
       coalesce (CASE age % 1 WHEN 0 THEN apoc.convert.toInteger(age) ELSE age END, '') AS Age,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(CASE weight % 1 WHEN 0 THEN apoc.convert.toInteger(weight) ELSE weight END, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`,
        coalesce(co.cohort_description, '') AS `Cohort`
'@

# "StudyFilesTab" row (row 4, column B): the case/parent-file query, with the
# RETURN clause reordered (Format moved up next to File Name, Size moved up
# next to File Type).
$studyFilesQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
MATCH (s:study)<-[*]-(c)<--(diag:diagnosis) ,(samp:sample)-->(c)
WHERE s.clinical_study_designation IN ['OSA01'] and demo.breed in ['Unknown'] and samp.summarized_sample_type IN ['Normal Cell Line']
WITH DISTINCT f,  s, c, demo, diag,parent,samp
WITH
        f, c, demo, diag, s,parent,samp,
        ['Bytes', 'KB', 'MB', 'GB', 'TB'] AS units,
        toInteger(floor(log(f.file_size)/log(1024))) as i,
        2 as precision
WITH    
        f, c, demo, diag, s,parent,samp,
        f.file_size /(1024^i) AS value, 10^precision AS factor,
        units[i] as unit
        WITH    
        f,  c, demo, diag, s, unit,parent,samp,
        round(factor * value)/factor AS size
RETURN DISTINCT
        coalesce(f.file_name, '') AS `File Name`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_type, '') AS `File Type`,
        CASE size % 1 WHEN 0 THEN apoc.convert.toInteger(size)+' ' +unit ELSE size+' ' +unit END AS Size,
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(samp.sample_id, '') AS `Sample ID`,
        coalesce(c.case_id, '') AS `Case ID`,
        coalesce(demo.breed,'') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis
'@

$rowHeight2 = $ws.Rows.Item(2).RowHeight
$rowHeight4 = $ws.Rows.Item(4).RowHeight

$ws.Range("B2").Value = $casesQuery
$ws.Range("B4").Value = $studyFilesQuery

# Re-apply the original row heights: they were explicit/custom, and should
# not change just because the cell text changed.
$ws.Rows.Item(2).RowHeight = $rowHeight2
$ws.Rows.Item(4).RowHeight = $rowHeight4

$null = $ws.Range("C4").Select()
